$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.333.56"
$ws.Range("E2").Value = "  -0.10%  "

$ws.Range("D3").Value = "'1.844.24"
$ws.Range("E3").Value = "  -1.06%  "

$ws.Range("D4").Value = "'0.9986"
$ws.Range("E4").Value = "  -0.20%  "

$ws.Range("D5").Value = "'233.02"
$ws.Range("E5").Value = "  -0.82%  "

$ws.Range("D6").Value = "'0.9991"
$ws.Range("E6").Value = "  -0.18%  "

$ws.Range("D7").Value = "'0.4646"
$ws.Range("E7").Value = "  -2.76%  "

$ws.Range("D8").Value = "'0.2733"
$ws.Range("E8").Value = "  -1.66%  "

$ws.Range("D9").Value = "'0.06284"
$ws.Range("E9").Value = "  -3.59%  "

$ws.Range("D10").Value = "'1.827.54"
$ws.Range("E10").Value = "  -1.84%  "

$ws.Range("D11").Value = "'0.07422"
$ws.Range("E11").Value = "  -0.29%  "

$ws.Range("D12").Value = "'16.28"
$ws.Range("E12").Value = "  +0.42%  "

$ws.Range("D13").Value = "'4.933"
$ws.Range("E13").Value = "  -2.12%  "

$ws.Range("D14").Value = "'83.81"
$ws.Range("E14").Value = "  -3.35%  "

$ws.Range("D15").Value = "'0.6216"
$ws.Range("E15").Value = "  -2.62%  "

$ws.Range("D16").Value = "'30.281.30"
$ws.Range("E16").Value = "  -0.18%  "

$ws.Range("D17").Value = "'0.9989"
$ws.Range("E17").Value = "  -0.15%  "

$ws.Range("D18").Value = "'228.47"
$ws.Range("E18").Value = "  -2.75%  "

$ws.Range("D19").Value = "'0.000007304"
$ws.Range("E19").Value = "  -1.80%  "

$ws.Range("E20").Value = "  -4.77%  "

$ws.Range("D21").Value = "'0.9984"
$ws.Range("E21").Value = "  -0.26%  "

$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").Value = "'4.924"
$ws.Range("E22").Value = "  -3.66%  "

$ws.Range("B23").Value = "Chainlink"
$ws.Range("C23").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D23").Value = "'5.864"
$ws.Range("E23").Value = "  -4.18%  "

$ws.Range("B24").Value = "Cosmos"
$ws.Range("C24").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D24").Value = "'9.195"
$ws.Range("E24").Value = "  -1.11%  "

$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D25").Value = "'164.67"
$ws.Range("E25").Value = "  -2.20%  "

$ws.Range("B26").Value = "EthereumClassic"
$ws.Range("C26").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D26").Value = "'17.80"
$ws.Range("E26").Value = "  -1.90%  "

$ws.Range("B27").Value = "LidoDAOToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D27").Value = "'1.871"
$ws.Range("E27").Value = "  -1.07%  "

$ws.Range("B28").Value = "Stellar"
$ws.Range("C28").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D28").Value = "'0.1032"
$ws.Range("E28").Value = "  -1.11%  "

$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "'1.371"
$ws.Range("E29").Value = "  -0.83%  "

$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D30").Value = "'4.085"
$ws.Range("E30").Value = "  -3.86%  "

$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").Value = "'3.814"
$ws.Range("E31").Value = "  -3.50%  "

$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").Value = "'0.04840"
$ws.Range("E32").Value = "  -2.51%  "

$ws.Range("B33").Value = "ARBITRUM"
$ws.Range("C33").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D33").Value = "'1.143"
$ws.Range("E33").Value = "  -1.80%  "

$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").Value = "'0.7106"
$ws.Range("E34").Value = "  -3.55%  "

$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").Value = "'2.697"
$ws.Range("E35").Value = "  -0.56%  "

$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D36").Value = "'0.01891"
$ws.Range("E36").Value = "  -2.03%  "

$ws.Range("B37").Value = "MXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D37").Value = "'2.659"
$ws.Range("E37").Value = "  +0.92%  "

$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").Value = "'0.8870"
$ws.Range("E38").Value = "  -2.74%  "

$ws.Range("B39").Value = "Quant"
$ws.Range("C39").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D39").Value = "'105.06"
$ws.Range("E39").Value = "  -0.57%  "

$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").Value = "'1.927"
$ws.Range("E40").Value = "  -4.56%  "

$ws.Range("B41").Value = "PaxDollar"
$ws.Range("C41").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D41").Value = "'1.002"
$ws.Range("E41").Value = "  +0.65%  "

$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").Value = "'5.561"
$ws.Range("E42").Value = "  -0.24%  "

$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D43").Value = "'0.4018"
$ws.Range("E43").Value = "  -3.51%  "

$ws.Range("B44").Value = "Aptos"
$ws.Range("C44").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D44").Value = "'7.081"
$ws.Range("E44").Value = "  -0.79%  "

$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").Value = "'60.93"
$ws.Range("E45").Value = "  -0.95%  "

$ws.Range("B46").Value = "Algorand"
$ws.Range("C46").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D46").Value = "'0.1197"
$ws.Range("E46").Value = "  -1.59%  "

$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "'8.574"
$ws.Range("E47").Value = "  -3.29%  "

$ws.Range("B48").Value = "Elrond"
$ws.Range("C48").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D48").Value = "'33.20"
$ws.Range("E48").Value = "  -0.54%  "

$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "'0.05505"
$ws.Range("E49").Value = "  -2.12%  "

$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").Value = "'1.351"
$ws.Range("E50").Value = "  -4.35%  "

$ws.Range("B51").Value = "Decentraland"
$ws.Range("C51").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D51").Value = "'0.3637"
$ws.Range("E51").Value = "  -3.15%  "
